$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Center-align the whole used range A2:D9
$ws.Range("A2:D9").HorizontalAlignment = -4108

# Row 6: replace the "X" placeholders in B6/C6 with real data
$ws.Range("B6").Value2 = 0.5
$ws.Range("C6").Formula = "= 10 / 60"

# Row 8: replace the "X" placeholders in B8/C8 with real data
$ws.Range("B8").Value2 = 0.5
$ws.Range("C8").Formula = "=10/60"

$ws.Calculate()
